$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 243, shifting the
# existing data (old rows 243-260) down to rows 245-262.
$ws.Rows("243:244").Insert()

# New row 243: Vega Monumental Concepción / Repollo / Crespo record / Primera
$ws.Range("A243").Value = 11
$ws.Range("B243").Value = "Vega Monumental Concepción"
$ws.Range("C243").Value = "Bíobío"
$ws.Range("D243").Value = 44568
$ws.Range("E243").Value = 8
$ws.Range("F243").Value = 100112006
$ws.Range("G243").Value = "Repollo"
$ws.Range("H243").Value = "Crespo record"
$ws.Range("I243").Value = "Primera"
$ws.Range("J243").Value = 1000
$ws.Range("K243").Value = 700
$ws.Range("L243").Value = 800
$ws.Range("M243").Value = 750
$ws.Range("N243").Value = "$/unidad"
$ws.Range("O243").Value = "Región Metropolitana"
$ws.Range("P243").Value = 750
$ws.Range("Q243").Value = 1
$ws.Range("R243").Value = "Hortaliza"

# New row 244: Vega Monumental Concepción / Repollo / Crespo record / Segunda
$ws.Range("A244").Value = 11
$ws.Range("B244").Value = "Vega Monumental Concepción"
$ws.Range("C244").Value = "Bíobío"
$ws.Range("D244").Value = 44568
$ws.Range("E244").Value = 8
$ws.Range("F244").Value = 100112006
$ws.Range("G244").Value = "Repollo"
$ws.Range("H244").Value = "Crespo record"
$ws.Range("I244").Value = "Segunda"
$ws.Range("J244").Value = 500
$ws.Range("K244").Value = 600
$ws.Range("L244").Value = 600
$ws.Range("M244").Value = 600
$ws.Range("N244").Value = "$/unidad"
$ws.Range("O244").Value = "Región Metropolitana"
$ws.Range("P244").Value = 600
$ws.Range("Q244").Value = 1
$ws.Range("R244").Value = "Hortaliza"
